$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menu Mock")
$ws.Rows.Item(22).Insert()
$v = $ws.Range("A23").Value()
Write-Host "A23 value after insert:"
Write-Host $v
$v2 = $ws.Range("A22").Value()
Write-Host "A22 value after insert:"
Write-Host $v2
